$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns
$ws.Range("C1").Value = "colecao"
$ws.Range("D1").Value = "vendedor"
$ws.Range("E1").Value = "regional"

# Copy the header style (bold + border) from A1/B1 onto the new headers
$ws.Range("A1:B1").Copy()
$ws.Range("C1:E1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Update existing row 2 data
$ws.Range("A2").Value = "Cliente A"
$ws.Range("B2").Value = 1000

# New row 3 data
$ws.Range("A3").Value = "Cliente B"
$ws.Range("B3").Value = 500

# New column data for rows 2 and 3
$ws.Range("C2").Value = "Coleção 1"
$ws.Range("C3").Value = "Coleção 2"
$ws.Range("D2").Value = "Vendedor 1"
$ws.Range("D3").Value = "Vendedor 2"
$ws.Range("E2").Value = "Região 1"
$ws.Range("E3").Value = "Região 2"
